# Added invivoPKfit outputs to dashboard script
# - Append a new "2.3.0" benchmark row to Table1 on Sheet1 (row 24)
# - Grow the table/autofilter range from A1:R23 to A1:R24
# - Record the new version string and notes text (these become new
#   shared-string entries automatically)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grow Table1 by one row (this also extends ref/autoFilter and the
# sheet's used range/dimension to A1:R24).
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add($tbl.ListRows.Count + 1)
$r = $newRow.Range

# Match the left-aligned style used by every other data row.
$r.HorizontalAlignment = -4131

$r.Cells.Item(1,1).Value  = "2.3.0"
$r.Cells.Item(1,2).Value  = 1023
$r.Cells.Item(1,3).Value  = 0.9999
$r.Cells.Item(1,4).Value  = 1
$r.Cells.Item(1,5).Value  = 1
$r.Cells.Item(1,6).Value  = 1.063
$r.Cells.Item(1,7).Value  = 352
$r.Cells.Item(1,8).Value  = 0.2996
$r.Cells.Item(1,9).Value  = 352
$r.Cells.Item(1,10).Value = 1.419
$r.Cells.Item(1,11).Value = 86
$r.Cells.Item(1,12).Value = 1.047
$r.Cells.Item(1,13).Value = 86
$r.Cells.Item(1,14).Value = 1.33
$r.Cells.Item(1,15).Value = 86
$r.Cells.Item(1,16).Value = 0.6344
$r.Cells.Item(1,17).Value = 863
$r.Cells.Item(1,18).Value = "Used Caco-2 to replace Fabs=Fgut=1"

# Leave the selection on the new row's last cell, like the author did
# after typing in the new data.
$ws.Range("R24").Select() | Out-Null
